$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1686.1282
$ws.Range("I15").Value = 1686.1282
$ws.Range("K15").Value = 5058.3846
$ws.Range("M15").Value = -4889.3846
$ws.Range("H17").Value = 2763.4546
$ws.Range("J17").Value = 2763.4546
$ws.Range("L17").Value = 8290.363799999999
$ws.Range("N17").Value = -8626.363799999999
$ws.Range("H33").Value = 481.66666
$ws.Range("I33").Value = 291.875
$ws.Range("K33").Value = 291.875
$ws.Range("M33").Value = -62.875
$ws.Range("H51").Value = 5371.4287
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 5371.4287
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = 5371.4287
$ws.Range("N51").Value = -6339.4287
$ws.Range("H63").Value = 39499.9
$ws.Range("I63").Value = 20000
$ws.Range("J63").Value = 117499.5
$ws.Range("K63").Value = 20000
$ws.Range("L63").Value = 117499.5
$ws.Range("M63").Value = -19376
$ws.Range("N63").Value = -118747.5
$ws.Range("H66").Value = 39499.9
$ws.Range("I66").Value = 20000
$ws.Range("J66").Value = 117499.5
$ws.Range("K66").Value = 60000
$ws.Range("L66").Value = 352498.5
$ws.Range("M66").Value = -56880
$ws.Range("N66").Value = -358738.5
$ws.Range("H68").Value = 76333
$ws.Range("J68").Value = 76333
$ws.Range("L68").Value = 76333
$ws.Range("N68").Value = -77831
$ws.Range("H71").Value = 76333
$ws.Range("J71").Value = 76333
$ws.Range("L71").Value = 228999
$ws.Range("N71").Value = -236487
$ws.Range("H92").Value = 4465387
$ws.Range("I92").Value = 862.5
$ws.Range("K92").Value = 862.5
$ws.Range("M92").Value = 385.5
$ws.Range("H100").Value = 4383.75
$ws.Range("I100").Value = 2247.0625
$ws.Range("J100").Value = 7232.6665
$ws.Range("K100").Value = 2247.0625
$ws.Range("L100").Value = 7232.6665
$ws.Range("M100").Value = -1706.0625
$ws.Range("N100").Value = -8314.666499999999
$ws.Range("H115").Value = 1435.8572
$ws.Range("I115").Value = 1408.6666
$ws.Range("K115").Value = 4225.9998
$ws.Range("M115").Value = -2658.9998
$ws.Range("H138").Value = 2999.0833
$ws.Range("J138").Value = 3016.5254
$ws.Range("L138").Value = 9049.5762
$ws.Range("N138").Value = -19329.5762

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3963.6667
$ws.Range("I45").Value = 2698.5
$ws.Range("J45").Value = 6494
$ws.Range("K45").Value = 2698.5
$ws.Range("L45").Value = 6494
$ws.Range("M45").Value = -2321.5
$ws.Range("N45").Value = -7248
$ws.Range("H61").Value = 4785.7207
$ws.Range("I61").Value = 3993.3928
$ws.Range("J61").Value = 6264.7334
$ws.Range("K61").Value = 3993.3928
$ws.Range("L61").Value = 6264.7334
$ws.Range("M61").Value = -3781.3928
$ws.Range("N61").Value = -6688.7334
$ws.Range("H97").Value = 877.3043
$ws.Range("I97").Value = 955.1905
$ws.Range("J97").Value = 59.5
$ws.Range("K97").Value = 955.1905
$ws.Range("L97").Value = 59.5
$ws.Range("M97").Value = -459.1905
$ws.Range("N97").Value = -1051.5
$ws.Range("H102").Value = 37038230
$ws.Range("I102").Value = 1344.875
$ws.Range("K102").Value = 1344.875
$ws.Range("M102").Value = 277.125
$ws.Range("H110").Value = 2165.375
$ws.Range("I110").Value = 2165.375
$ws.Range("K110").Value = 2165.375
$ws.Range("M110").Value = -120.375
$ws.Range("H132").Value = 3431.0444
$ws.Range("I132").Value = 3395.3865
$ws.Range("K132").Value = 10186.1595
$ws.Range("M132").Value = -7656.1595
$ws.Range("H136").Value = 4785.7207
$ws.Range("I136").Value = 3993.3928
$ws.Range("J136").Value = 6264.7334
$ws.Range("K136").Value = 11980.1784
$ws.Range("L136").Value = 18794.2002
$ws.Range("M136").Value = -9430.178400000001
$ws.Range("N136").Value = -23894.2002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 50335.24
$ws.Range("I20").Value = 2866.3845
$ws.Range("J20").Value = 127472.125
$ws.Range("K20").Value = 2866.3845
$ws.Range("L20").Value = 127472.125
$ws.Range("M20").Value = -2619.3845
$ws.Range("N20").Value = -127966.125
$ws.Range("H107").Value = 2807.2
$ws.Range("I107").Value = 2505.5
$ws.Range("J107").Value = 3008.3333
$ws.Range("K107").Value = 2505.5
$ws.Range("L107").Value = 3008.3333
$ws.Range("M107").Value = -585.5
$ws.Range("N107").Value = -6848.3333
$ws.Range("H134").Value = 3848.889
$ws.Range("I134").Value = 2460.3
$ws.Range("J134").Value = 7816.2856
$ws.Range("K134").Value = 7380.900000000001
$ws.Range("L134").Value = 23448.8568
$ws.Range("M134").Value = -4845.900000000001
$ws.Range("N134").Value = -28518.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 487.25
$ws.Range("J15").Value = 487.25
$ws.Range("L15").Value = 487.25
$ws.Range("N15").Value = -827.25
$ws.Range("H31").Value = 3511.2646
$ws.Range("I31").Value = 2637.2173
$ws.Range("J31").Value = 5338.8184
$ws.Range("K31").Value = 2637.2173
$ws.Range("L31").Value = 5338.8184
$ws.Range("M31").Value = -2342.2173
$ws.Range("N31").Value = -5928.8184
$ws.Range("H34").Value = 3511.2646
$ws.Range("I34").Value = 2637.2173
$ws.Range("J34").Value = 5338.8184
$ws.Range("K34").Value = 2637.2173
$ws.Range("L34").Value = 5338.8184
$ws.Range("M34").Value = -2435.2173
$ws.Range("N34").Value = -5742.8184
$ws.Range("H68").Value = 40295
$ws.Range("J68").Value = 40295
$ws.Range("L68").Value = 40295
$ws.Range("N68").Value = -41793
$ws.Range("H71").Value = 40295
$ws.Range("J71").Value = 40295
$ws.Range("L71").Value = 120885
$ws.Range("N71").Value = -128373
$ws.Range("H106").Value = 94999.5
$ws.Range("J106").Value = 94999.5
$ws.Range("L106").Value = 94999.5
$ws.Range("N106").Value = -97523.5
$ws.Range("H132").Value = 3264.5
$ws.Range("I132").Value = 3264.5
$ws.Range("K132").Value = 9793.5
$ws.Range("M132").Value = -7263.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 14399.6
$ws.Range("J55").Value = 14399.6
$ws.Range("L55").Value = 14399.6
$ws.Range("N55").Value = -15053.6
$ws.Range("H70").Value = 55615.727
$ws.Range("I70").Value = 128489.445
$ws.Range("J70").Value = 5164.6924
$ws.Range("K70").Value = 128489.445
$ws.Range("L70").Value = 5164.6924
$ws.Range("M70").Value = -128219.445
$ws.Range("N70").Value = -5704.6924
$ws.Range("H73").Value = 55615.727
$ws.Range("I73").Value = 128489.445
$ws.Range("J73").Value = 5164.6924
$ws.Range("K73").Value = 128489.445
$ws.Range("L73").Value = 5164.6924
$ws.Range("M73").Value = -127553.445
$ws.Range("N73").Value = -7036.6924
$ws.Range("H136").Value = 18737.125
$ws.Range("J136").Value = 18737.125
$ws.Range("L136").Value = 56211.375
$ws.Range("N136").Value = -61311.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4859
$ws.Range("I7").Value = 2502.6667
$ws.Range("K7").Value = 2502.6667
$ws.Range("M7").Value = -2390.6667
$ws.Range("H40").Value = 12768.333
$ws.Range("I40").Value = 14479.444
$ws.Range("J40").Value = 7635
$ws.Range("K40").Value = 14479.444
$ws.Range("L40").Value = 7635
$ws.Range("M40").Value = -14343.444
$ws.Range("N40").Value = -7907
$ws.Range("H126").Value = 4859
$ws.Range("I126").Value = 2502.6667
$ws.Range("K126").Value = 7508.000100000001
$ws.Range("M126").Value = -5038.000100000001
$ws.Range("H132").Value = 4493.037
$ws.Range("I132").Value = 3830.9565
$ws.Range("J132").Value = 8300
$ws.Range("K132").Value = 11492.8695
$ws.Range("L132").Value = 24900
$ws.Range("M132").Value = -8962.869499999999
$ws.Range("N132").Value = -29960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 685.55
$ws.Range("I100").Value = 717.4706
$ws.Range("J100").Value = 504.66666
$ws.Range("K100").Value = 1434.9412
$ws.Range("L100").Value = 1009.33332
$ws.Range("M100").Value = -893.9412
$ws.Range("N100").Value = -2091.33332
$ws.Range("H122").Value = 2874.4119
$ws.Range("I122").Value = 2091
$ws.Range("J122").Value = 8750
$ws.Range("K122").Value = 6273
$ws.Range("L122").Value = 26250
$ws.Range("M122").Value = -3823
$ws.Range("N122").Value = -31150
$ws.Range("H126").Value = 1576.4117
$ws.Range("I126").Value = 1386
$ws.Range("J126").Value = 2195.25
$ws.Range("K126").Value = 4158
$ws.Range("L126").Value = 6585.75
$ws.Range("M126").Value = -1688
$ws.Range("N126").Value = -11525.75
$ws.Range("H132").Value = 1560.3158
$ws.Range("I132").Value = 1202.5555
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 3607.6665
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -1077.6665
$ws.Range("N132").Value = -29060
$ws.Range("H136").Value = 2997.3044
$ws.Range("I136").Value = 2217.7297
$ws.Range("K136").Value = 6653.1891
$ws.Range("M136").Value = -4103.1891
